$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 803, shifting existing rows 803-861 down to 804-862.
$ws.Rows.Item(803).Insert()

# Populate the newly inserted row 803 with the new weekly record.
$ws.Cells.Item(803, 1).Value2 = 3
$ws.Cells.Item(803, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(803, 3).Value2 = "Coquimbo"
$ws.Cells.Item(803, 4).Value2 = 45265
$ws.Cells.Item(803, 5).Value2 = 5
$ws.Cells.Item(803, 6).Value2 = 100112037
$ws.Cells.Item(803, 7).Value2 = "Cebollín"
$ws.Cells.Item(803, 8).Value2 = "Sin especificar"
$ws.Cells.Item(803, 9).Value2 = "Primera"
$ws.Cells.Item(803, 10).Value2 = 150
$ws.Cells.Item(803, 11).Value2 = 4000
$ws.Cells.Item(803, 12).Value2 = 4000
$ws.Cells.Item(803, 13).Value2 = 4000
$ws.Cells.Item(803, 14).Value2 = "$/paquete 36 unidades"
$ws.Cells.Item(803, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(803, 16).Value2 = 111
$ws.Cells.Item(803, 17).Value2 = 36
$ws.Cells.Item(803, 18).Value2 = "Hortaliza"
